$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.254.79"
$ws.Range("E2").Value = "  -1.54%  "

$ws.Range("D3").Value = "2.574.25"
$ws.Range("E3").Value = "  -1.89%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'555.62"
$ws.Range("E5").Value = "  -2.01%  "

$ws.Range("D6").Value = "'141.61"
$ws.Range("E6").Value = "  -2.61%  "

$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("E8").Value = "  -1.58%  "

$ws.Range("D9").Value = "2.581.23"
$ws.Range("E9").Value = "  -2.41%  "

$ws.Range("E10").Value = "  -1.51%  "

$ws.Range("E11").Value = "  -0.97%  "

$ws.Range("D12").Value = "'0.164"
$ws.Range("E12").Value = "  +10.95%  "

$ws.Range("E13").Value = "  +2.31%  "

$ws.Range("D14").Value = "3.032.30"
$ws.Range("E14").Value = "  -2.01%  "

$ws.Range("D15").Value = "59.255.54"
$ws.Range("E15").Value = "  -1.50%  "

$ws.Range("D16").Value = "'22.93"
$ws.Range("E16").Value = "  +4.41%  "

$ws.Range("E17").Value = "  -0.38%  "

$ws.Range("D18").Value = "2.578.06"
$ws.Range("E18").Value = "  -2.52%  "

$ws.Range("E19").Value = "  +0.34%  "

$ws.Range("D20").Value = "'337.98"
$ws.Range("E20").Value = "  -1.10%  "

$ws.Range("E21").Value = "  -0.72%  "

$ws.Range("E22").Value = "  +1.30%  "

$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("E24").Value = "  +9.48%  "

$ws.Range("D25").Value = "'62.60"
$ws.Range("E25").Value = "  -4.89%  "

$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("E27").Value = "  -2.65%  "

$ws.Range("D28").Value = "'7.37"
$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").Value = "0.0₃0772"
$ws.Range("E29").Value = "  -3.50%  "

$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("E31").Value = "  +0.91%  "

$ws.Range("D32").Value = "'1.67"
$ws.Range("E32").Value = "  -1.60%  "

$ws.Range("D33").Value = "'159.34"
$ws.Range("E33").Value = "  +0.28%  "

$ws.Range("E34").Value = "  -0.39%  "

$ws.Range("E35").Value = "  -0.30%  "

$ws.Range("E36").Value = "  +1.33%  "

$ws.Range("D37").Value = "'0.894"
$ws.Range("E37").Value = "  +1.16%  "

$ws.Range("E38").Value = "  -0.53%  "

$ws.Range("E39").Value = "  -4.00%  "

$ws.Range("E40").Value = "  -1.89%  "

$ws.Range("E41").Value = "  +1.11%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'138.42"
$ws.Range("E42").Value = "  +8.58%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'288.99"
$ws.Range("E43").Value = "  -3.43%  "

$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.42%  "

$ws.Range("E45").Value = "  -1.09%  "

$ws.Range("E46").Value = "  -1.73%  "

$ws.Range("E47").Value = "  -0.11%  "

$ws.Range("E48").Value = "  -2.72%  "

$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("D50").Value = "'18.62"
$ws.Range("E50").Value = "  -0.27%  "

$ws.Range("D51").Value = "1.938.91"
$ws.Range("E51").Value = "  -1.02%  "
